# Rename the header cell from "Name" to "Tag" and refresh the formatting
# applied to the header/data cells, then move the active selection to A2
# (matching the frozen-pane header row selection recorded in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 held the column header "Name" - rename it to "Tag".
$ws.Range("A1").Value = "Tag"

# Re-apply the (Normal) cell style to the header + data cell so the sheet
# carries an explicit formatting record for both, matching the refreshed
# style that ships with the renamed column.
$ws.Range("A1:A2").Style = "Normal"

# Move/record the active selection on the frozen pane to A2.
$ws.Range("A2").Select()
